$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46056, 5.63, 0.037),
    @(46056.01041666666, 3.721, 1.744),
    @(46056.02083333334, 5.369, 0.324),
    @(46056.03125, 17.251, 0),
    @(46056.04166666666, 3.608, 3.09),
    @(46056.05208333334, 5.79, 0.618),
    @(46056.0625, 12.878, 0),
    @(46056.07291666666, 2.363, 0.318),
    @(46056.08333333334, 18.51, 0),
    @(46056.09375, 9.106, 0.484),
    @(46056.10416666666, 4.966, 0.114),
    @(46056.11458333334, 12.135, 0),
    @(46056.125, 7.856, 0.926),
    @(46056.13541666666, 9.939, 0),
    @(46056.14583333334, 15.076, 0),
    @(46056.15625, 6.875, 0.112),
    @(46056.16666666666, 15.514, 0.008),
    @(46056.17708333334, 23.624, 0),
    @(46056.1875, 7.982, 0),
    @(46056.19791666666, 24.11, 0),
    @(46056.20833333334, 1.651, 5.887),
    @(46056.21875, 6.321, 0.027),
    @(46056.22916666666, 0.636, 2.821),
    @(46056.23958333334, 8.433, 0),
    @(46056.25, 0.106, 7.406),
    @(46056.26041666666, 3.375, 0),
    @(46056.27083333334, 7.466, 0),
    @(46056.28125, 0.6850000000000001, 0),
    @(46056.29166666666, 6.824, 0),
    @(46056.30208333334, 0.6919999999999999, 1.071),
    @(46056.3125, 0.074, 1.898),
    @(46056.32291666666, 0.616, 3.514),
    @(46056.33333333334, 0, 19.937),
    @(46056.34375, 0, 12.123),
    @(46056.35416666666, 2.026, 0.07099999999999999),
    @(46056.36458333334, 0.131, 7.083),
    @(46056.375, 4.382, 0.708),
    @(46056.38541666666, 0.24, 8.720000000000001),
    @(46056.39583333334, 0, 12.119),
    @(46056.40625, 0.002, 17.915),
    @(46056.41666666666, 0.653, 4.153),
    @(46056.42708333334, 1.402, 0),
    @(46056.4375, 0.158, 19.516),
    @(46056.44791666666, 0, 30.03),
    @(46056.45833333334, 4.261, 2.079),
    @(46056.46875, 1.283, 3.205),
    @(46056.47916666666, 1.261, 1.158),
    @(46056.48958333334, 2.389, 0.5669999999999999),
    @(46056.5, 1.821, 0),
    @(46056.51041666666, 0.863, 4.529),
    @(46056.53125, 0, 0),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# New rows 48-52 need the same date/time number format as column A above them
$ws.Range("A48:A52").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A1").Select() | Out-Null
